$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3569133333333334
$ws.Range("H2").Value = 1.07074
$ws.Range("I2").Value = 0.928094236140756
$ws.Range("J2").Value = 0.9508856218001945
$ws.Range("M2").Value = 6.904282333333334
$ws.Range("N2").Value = 20.712847
$ws.Range("O2").Value = 0.04090096694673821
$ws.Range("P2").Value = 0.04187630379952963
$ws.Range("Q2").Value = 2.464230421864445
$ws.Range("R2").Value = 22.17807379678
$ws.Range("S2").Value = 0.0379599516758513
$ws.Range("T2").Value = 0.03981957517710958
$ws.Range("G3").Value = 0.3569133333333334
$ws.Range("H3").Value = 1.07074
$ws.Range("I3").Value = 0.928094236140756
$ws.Range("J3").Value = 0.9508856218001945
$ws.Range("O3").Value = 0.08154120089165004
$ws.Range("P3").Value = 0.08348565707905657
$ws.Range("Q3").Value = 4.91275201718889
$ws.Range("R3").Value = 44.2147681547
$ws.Range("S3").Value = 0.07567791855553588
$ws.Range("T3").Value = 0.07938531094301653
$ws.Range("G4").Value = 0.3569133333333334
$ws.Range("H4").Value = 1.07074
$ws.Range("I4").Value = 0.928094236140756
$ws.Range("J4").Value = 0.9508856218001945
$ws.Range("M4").Value = 68.52477533333332
$ws.Range("N4").Value = 205.574326
$ws.Range("O4").Value = 0.4059407532351291
$ws.Range("P4").Value = 0.4156209394565383
$ws.Range("Q4").Value = 24.45740598013778
$ws.Range("R4").Value = 220.11665382124
$ws.Range("S4").Value = 0.3767512732921602
$ws.Range("T4").Value = 0.3952079754483114
$ws.Range("G5").Value = 0.3569133333333334
$ws.Range("H5").Value = 1.07074
$ws.Range("I5").Value = 0.928094236140756
$ws.Range("J5").Value = 0.9508856218001945
$ws.Range("M5").Value = 11.7948525
$ws.Range("N5").Value = 23.589705
$ws.Range("O5").Value = 0.0698727034836137
$ws.Range("P5").Value = 0.04769260609713784
$ws.Range("Q5").Value = 4.20974012195
$ws.Range("R5").Value = 25.2584407317
$ws.Range("S5").Value = 0.06484845336671399
$ws.Range("T5").Value = 0.04535021340394866
$ws.Range("G6").Value = 0.3569133333333334
$ws.Range("H6").Value = 1.07074
$ws.Range("I6").Value = 0.928094236140756
$ws.Range("J6").Value = 0.9508856218001945
$ws.Range("M6").Value = 67.81640633333333
$ws.Range("N6").Value = 203.449219
$ws.Range("O6").Value = 0.401744375442869
$ws.Range("P6").Value = 0.4113244935677377
$ws.Range("Q6").Value = 24.20457963911778
$ws.Range("R6").Value = 217.84121675206
$ws.Range("S6").Value = 0.3728566392504946
$ws.Range("T6").Value = 0.3911225468278084
$ws.Range("G7").Value = 0.0276525
$ws.Range("H7").Value = 0.055305
$ws.Range("I7").Value = 0.07190576385924385
$ws.Range("J7").Value = 0.04911437819980551
$ws.Range("M7").Value = 6.904282333333334
$ws.Range("N7").Value = 20.712847
$ws.Range("O7").Value = 0.04090096694673821
$ws.Range("P7").Value = 0.04187630379952963
$ws.Range("Q7").Value = 0.1909206672225
$ws.Range("R7").Value = 1.145524003335
$ws.Range("S7").Value = 0.002941015270886895
$ws.Range("T7").Value = 0.002056728622420051
$ws.Range("G8").Value = 0.0276525
$ws.Range("H8").Value = 0.055305
$ws.Range("I8").Value = 0.07190576385924385
$ws.Range("J8").Value = 0.04911437819980551
$ws.Range("Q8").Value = 0.3806242649625
$ws.Range("R8").Value = 2.283745589775
$ws.Range("S8").Value = 0.005863282336114152
$ws.Range("T8").Value = 0.004100346136040055
$ws.Range("G9").Value = 0.0276525
$ws.Range("H9").Value = 0.055305
$ws.Range("I9").Value = 0.07190576385924385
$ws.Range("J9").Value = 0.04911437819980551
$ws.Range("M9").Value = 68.52477533333332
$ws.Range("N9").Value = 205.574326
$ws.Range("O9").Value = 0.4059407532351291
$ws.Range("P9").Value = 0.4156209394565383
$ws.Range("Q9").Value = 1.894881349905
$ws.Range("R9").Value = 11.36928809943
$ws.Range("S9").Value = 0.02918947994296877
$ws.Range("T9").Value = 0.02041296400822689
$ws.Range("G10").Value = 0.0276525
$ws.Range("H10").Value = 0.055305
$ws.Range("I10").Value = 0.07190576385924385
$ws.Range("J10").Value = 0.04911437819980551
$ws.Range("M10").Value = 11.7948525
$ws.Range("N10").Value = 23.589705
$ws.Range("O10").Value = 0.0698727034836137
$ws.Range("P10").Value = 0.04769260609713784
$ws.Range("Q10").Value = 0.32615715875625
$ws.Range("R10").Value = 1.304628635025
$ws.Range("S10").Value = 0.005024250116899692
$ws.Range("T10").Value = 0.002342392693189178
$ws.Range("G11").Value = 0.0276525
$ws.Range("H11").Value = 0.055305
$ws.Range("I11").Value = 0.07190576385924385
$ws.Range("J11").Value = 0.04911437819980551
$ws.Range("M11").Value = 67.81640633333333
$ws.Range("N11").Value = 203.449219
$ws.Range("O11").Value = 0.401744375442869
$ws.Range("P11").Value = 0.4113244935677377
$ws.Range("Q11").Value = 1.8752931761325
$ws.Range("R11").Value = 11.251759056795
$ws.Range("S11").Value = 0.02888773619237434
$ws.Range("T11").Value = 0.02020194673992934